$wb = $excel.ActiveWorkbook

# ALC row 17: One for the Road
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 974.86664
$ws.Range("J17").Value = 849.88
$ws.Range("L17").Value = 2549.64
$ws.Range("N17").Value = -2885.64

# ALC row 62: The Mustache Suits Him
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5901.769
$ws.Range("I62").Value = 5178.6
$ws.Range("J62").Value = 8312.333000000001
$ws.Range("K62").Value = 5178.6
$ws.Range("L62").Value = 8312.333000000001
$ws.Range("M62").Value = -4554.6
$ws.Range("N62").Value = -9560.333000000001

# ALC row 65: Forgery of Convenience (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5901.769
$ws.Range("I65").Value = 5178.6
$ws.Range("J65").Value = 8312.333000000001
$ws.Range("K65").Value = 25893
$ws.Range("L65").Value = 41561.665
$ws.Range("M65").Value = -22773
$ws.Range("N65").Value = -47801.665

# ALC row 129: Practical Command
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1425.5
$ws.Range("I129").Value = 458.92307
$ws.Range("K129").Value = 1376.76921
$ws.Range("M129").Value = 3623.23079

# ALC row 132: Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5383.838
$ws.Range("I132").Value = 4648.5806
$ws.Range("K132").Value = 13945.7418
$ws.Range("M132").Value = -11415.7418

# ALC row 135: For Tired Minds
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1589.2858
$ws.Range("I135").Value = 225.6
$ws.Range("K135").Value = 2030.4
$ws.Range("M135").Value = 504.6000000000001

# ALC row 138: All-night Crafting
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2859.6086
$ws.Range("I138").Value = 675.38464
$ws.Range("J138").Value = 5699.1
$ws.Range("K138").Value = 2026.15392
$ws.Range("L138").Value = 17097.3
$ws.Range("M138").Value = 3113.84608
$ws.Range("N138").Value = -27377.3

# ARM row 63: Rivets Run through It
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()

# ARM row 66: A Riveting Revival (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()

# ARM row 110: Scheduled Maintenance
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1099.7693
$ws.Range("I110").Value = 921.8333
$ws.Range("K110").Value = 921.8333
$ws.Range("M110").Value = 1123.1667

# ARM row 122: Haste for High Durium
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5001.7144
$ws.Range("I122").Value = 5001.7144
$ws.Range("K122").Value = 15005.1432
$ws.Range("M122").Value = -12555.1432

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2694.05
$ws.Range("I132").Value = 2752.818
$ws.Range("J132").Value = 2622.2222
$ws.Range("K132").Value = 8258.454000000002
$ws.Range("L132").Value = 7866.6666
$ws.Range("M132").Value = -5728.454000000002
$ws.Range("N132").Value = -12926.6666

# ARM row 135: Forgiveness for My Shins
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# BSM row 92: Have Blade, Will Travel
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 32249
$ws.Range("J92").Value = 32249
$ws.Range("L92").Value = 32249
$ws.Range("N92").Value = -37241

# BSM row 134: Ruthenium Supremium
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3113.75
$ws.Range("I134").Value = 2582
$ws.Range("K134").Value = 7746
$ws.Range("M134").Value = -5211

# CRP row 31: Wall Not Found
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3894.182
$ws.Range("I31").Value = 3080.9285
$ws.Range("K31").Value = 3080.9285
$ws.Range("M31").Value = -2785.9285

# CRP row 34: Armoires of the Rich and Famous
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3894.182
$ws.Range("I34").Value = 3080.9285
$ws.Range("K34").Value = 3080.9285
$ws.Range("M34").Value = -2878.9285

# CRP row 59: Bow Down to Magic
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 78563.5
$ws.Range("J59").Value = 78563.5
$ws.Range("L59").Value = 78563.5
$ws.Range("N59").Value = -80853.5

# CRP row 60: Bowing to Greater Power
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 18174.834
$ws.Range("I60").Value = 14000
$ws.Range("K60").Value = 14000
$ws.Range("M60").Value = -13489

# CRP row 62: Splinter in the Sewers
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# CRP row 65: The Lumber of Their Discontent (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# CRP row 68: Do You Even String Bow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# CRP row 71: Win One Bow, Get Three Free (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# CRP row 74: License to Heal
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# CRP row 77: Purified Polyrhythm (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# CRP row 94: Beech, Please
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 631.7
$ws.Range("J94").Value = 782.6
$ws.Range("L94").Value = 782.6
$ws.Range("N94").Value = -1684.6

# CRP row 132: Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5356
$ws.Range("I132").Value = 4651.7334
$ws.Range("J132").Value = 7116.6665
$ws.Range("K132").Value = 13955.2002
$ws.Range("L132").Value = 21349.9995
$ws.Range("M132").Value = -11425.2002
$ws.Range("N132").Value = -26409.9995

# CUL row 34: Fever Pitch
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 365
$ws.Range("I34").Value = 365
$ws.Range("K34").Value = 1095
$ws.Range("M34").Value = -1011

# CUL row 39: Bloody Good Tart, This
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4449
$ws.Range("J39").Value = 6398
$ws.Range("L39").Value = 19194
$ws.Range("N39").Value = -19782

# CUL row 56: Culture Club
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5533.2856
$ws.Range("I56").Value = 5533.2856
$ws.Range("K56").Value = 5533.2856
$ws.Range("M56").Value = -5003.2856

# GSM row 113: Copious Crystal Cannons
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2187.5
$ws.Range("I113").Value = 2231
$ws.Range("K113").Value = 2231
$ws.Range("M113").Value = -61

# GSM row 122: Awarding Academic Excellence
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2928.0356
$ws.Range("I122").Value = 2019.2142
$ws.Range("J122").Value = 3836.8572
$ws.Range("K122").Value = 6057.642599999999
$ws.Range("L122").Value = 11510.5716
$ws.Range("M122").Value = -3607.642599999999
$ws.Range("N122").Value = -16410.5716

# GSM row 132: On Board for Lar
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2076.2307
$ws.Range("I132").Value = 1739.1
$ws.Range("K132").Value = 5217.299999999999
$ws.Range("M132").Value = -2687.299999999999

# LTW row 16: Saddle Sore
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2340

# LTW row 40: Best Served Toad
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 49000
$ws.Range("I40").Value = 49000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 49000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -48864
$ws.Range("N40").ClearContents()

# LTW row 61: Spelling Me Softly
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

# LTW row 113: Peace in Rest
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# LTW row 122: Hell on Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7942.6665
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 7942.6665
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 23827.9995
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -28727.9995

# LTW row 136: Respect for Br'aax
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6812.0713
$ws.Range("I136").Value = 6767.5713
$ws.Range("K136").Value = 20302.7139
$ws.Range("M136").Value = -17752.7139

# WVR row 94: Proper Props
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# WVR row 113: A Tender Table
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 536.625
$ws.Range("I113").Value = 533.8333
$ws.Range("J113").Value = 545
$ws.Range("K113").Value = 1601.4999
$ws.Range("L113").Value = 1635
$ws.Range("M113").Value = 568.5001
$ws.Range("N113").Value = -5975

# WVR row 126: A Polished Purchase
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4299
$ws.Range("I126").Value = 4299
$ws.Range("K126").Value = 12897
$ws.Range("M126").Value = -10427

# WVR row 136: Weaving the Envelope
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4059.3076
$ws.Range("I136").Value = 4161
$ws.Range("K136").Value = 12483
$ws.Range("M136").Value = -9933
